$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: update responsible name ---
$ws.Range("C4").Value = "Marisol Ornelas"

# --- Row 5 ---
$ws.Range("B5").Value = "No todos los proyectos cuentan con carta de aceptación"
$ws.Range("C5").Value = "Equipo de ventas"
$ws.Range("D5").NumberFormat = "DD/MM/YY"
$ws.Range("D5").Value = (Get-Date -Year 2016 -Month 2 -Day 4).Date
$ws.Range("F5").Value = "En proceso"
$ws.Range("G5").Value = "Se solicita validar que se esta enviando la carta para cada ticket resuelto"
$ws.Range("A5").RowHeight = 41.75

# --- Row 6 ---
$ws.Range("B6").Value = "No se tiene realizadas encuestas de satisfacción"
$ws.Range("C6").Value = "Magda Montoya"
$ws.Range("D6").NumberFormat = "DD/MM/YY"
$ws.Range("D6").Value = (Get-Date -Year 2016 -Month 2 -Day 4).Date
$ws.Range("F6").Value = "En proceso"
$ws.Range("G6").Value = "Realizar encuesta de satisfacción a los proyectos señalados"
$ws.Range("A6").RowHeight = 41.75

# --- Row 7 ---
$ws.Range("B7").Value = "No se tiene establecido el nombre adecuado en los archivos de salida por el proceso de ventas"
$ws.Range("C7").Value = "Equipo de ventas"
$ws.Range("D7").NumberFormat = "DD/MM/YY"
$ws.Range("D7").Value = (Get-Date -Year 2016 -Month 2 -Day 3).Date
$ws.Range("F7").Value = "En proceso"
$ws.Range("G7").Value = "apegarse al nombre de los archivos en el plan de configuración"
$ws.Range("A7").RowHeight = 41.75

# --- Row 8 ---
$ws.Range("B8").Value = "No se esta respetando las ubicación física de los archivos establecido por el plan de configuración"
$ws.Range("C8").Value = "Equipo de ventas"
$ws.Range("D8").NumberFormat = "DD/MM/YY"
$ws.Range("D8").Value = (Get-Date -Year 2016 -Month 2 -Day 3).Date
$ws.Range("F8").Value = "En proceso"
$ws.Range("G8").Value = "Apegarse a la dirección establecida en el plan de configuración"
$ws.Range("A8").RowHeight = 41.75

# --- Update active selection to A8 ---
$ws.Range("A8").Select()
